# Ind_Int_Mex.xlsx edit script
# - Fix "Anoi" typo to "Ano" in the table header (cell B4)
# - Change the Month column (C5:C84) from numeric month values (1-12)
#   to abbreviated Spanish month-name text (Ene., Feb., Mar., ...)
# - Replace the footnote in B89 with the shorter-indent version and
#   apply a left-aligned, indented cell style to it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix header typo: "Anoi" -> "Ano" ---
$ws.Range("B4").Value2 = "A" + [char]0x00F1 + "o"

# --- 2. Convert month numbers to abbreviated month names ---
$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

for ($r = 5; $r -le 84; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $num = [int]$cell.Value2
    $cell.Value2 = $monthNames[$num]
}

# Normalize row heights for the rows we touched (writing into a hidden
# row recalculates an autofit height as a side effect); restore the
# rows that genuinely carry an explicit custom height afterwards.
for ($r = 5; $r -le 84; $r++) {
    $ws.Rows.Item($r).AutoFit()
}
$ws.Rows.Item(5).RowHeight = 18
$ws.Rows.Item(6).RowHeight = 18
$ws.Rows.Item(37).RowHeight = 18
$ws.Rows.Item(83).RowHeight = 18

# --- 3. Update the footnote cell B89 ---
$paraText = "   Para variaci" + [char]0x00F3 + "n porcentual anual: C" + [char]0x00E1 + "lculos propios con base en datos proporcionados por el INEGI."
$ws.Range("B89").Value2 = $paraText
$ws.Range("B89").HorizontalAlignment = -4131
$ws.Range("B89").IndentLevel = 4

Write-Host "Edit complete"
